$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the current "accounts" sheet -> becomes "noAccounts".
#    The duplicate is placed immediately after "accounts" and keeps
#    the ORIGINAL (old) "accounts" data/layout untouched (A1:H2),
#    only its selection is updated.
# ------------------------------------------------------------------
$accounts = $wb.Worksheets.Item("accounts")
$accounts.Copy($null, $accounts)
$noAccounts = $wb.Worksheets.Item("accounts (2)")
$noAccounts.Name = "noAccounts"
$noAccounts.Range("E4").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Replace the content of the (original) "accounts" sheet with the
#    new layout: data now starts at column B (column A left blank)
#    and two new trailing columns are added (textOffer / txtAccount)
#    plus a brand-new "txtBeneficiary" column.
# ------------------------------------------------------------------
$accounts.Cells.Clear()

$accounts.Range("B1").Value = "textOnboarding"
$accounts.Range("C1").Value = "textLogin"
$accounts.Range("D1").Value = "idType"
$accounts.Range("E1").Value = "idNumber"
$accounts.Range("F1").Value = "textpassword"
$accounts.Range("G1").Value = "Password"
$accounts.Range("H1").Value = "loader"
$accounts.Range("I1").Value = "textOffer"
$accounts.Range("J1").Value = "txtAccount"
$accounts.Range("K1").Value = "txtBeneficiary"

$accounts.Range("B2").Value = "Le damos la bienvenida a"
$accounts.Range("C2").Value = "Bienvenido a"
$accounts.Range("D2").Value = "Cédula de ciudadanía"
$accounts.Range("E2").Value = 1013583153
$accounts.Range("F2").Value = "Por su seguridad"
$accounts.Range("G2").Value = "ibcs0011"
$accounts.Range("H2").Value = "validando"
$accounts.Range("I2").Value = "Personalice su oferta"
$accounts.Range("J2").Value = "Seleccionar cuenta"
$accounts.Range("K2").Value = "Seguro de vida"

# left-align style to mirror the header/id columns used elsewhere
$accounts.Range("C1").HorizontalAlignment = -4131
$accounts.Range("E1").HorizontalAlignment = -4131
$accounts.Range("E2").HorizontalAlignment = -4131
$accounts.Range("F2").HorizontalAlignment = -4131

$accounts.Range("K3").Select() | Out-Null
$accounts.Activate() | Out-Null

# ------------------------------------------------------------------
# 3. "offer" sheet: the account id used for the demo row changes.
# ------------------------------------------------------------------
$offer = $wb.Worksheets.Item("offer")
$offer.Range("D2").Value = 1013583153
$offer.Range("A1:I2").Select() | Out-Null

# ------------------------------------------------------------------
# 4. "consolidated" sheet: no data change, the tab just stops being
#    the active one (handled by activating "accounts" above).
# ------------------------------------------------------------------
$consolidated = $wb.Worksheets.Item("consolidated")
$consolidated.Range("C8").Select() | Out-Null

# Re-activate "accounts" so it ends up as the active sheet/tab.
$accounts.Activate() | Out-Null
$accounts.Range("K3").Select() | Out-Null
